# Update TPM-derived statistics in the LR-pairs sheet with newly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.999699333333333
$ws.Range("H2").Value = 5.999098
$ws.Range("I2").Value = 0.5605459117818491
$ws.Range("J2").Value = 0.5605459117818491
$ws.Range("M2").Value = 0.08962966666666666
$ws.Range("N2").Value = 0.268889
$ws.Range("O2").Value = 0.4339761198462219
$ws.Range("P2").Value = 0.4339761198462219
$ws.Range("Q2").Value = 0.1792323846802222
$ws.Range("R2").Value = 1.613091462122
$ws.Range("S2").Value = 0.2432635397907494
$ws.Range("T2").Value = 0.2432635397907494

# Row 3
$ws.Range("G3").Value = 1.999699333333333
$ws.Range("H3").Value = 5.999098
$ws.Range("I3").Value = 0.5605459117818491
$ws.Range("J3").Value = 0.5605459117818491
$ws.Range("M3").Value = 0.1169016666666667
$ws.Range("N3").Value = 0.350705
$ws.Range("O3").Value = 0.5660238801537781
$ws.Range("P3").Value = 0.5660238801537781
$ws.Range("Q3").Value = 0.2337681848988889
$ws.Range("R3").Value = 2.10391366409
$ws.Range("S3").Value = 0.3172823719910997
$ws.Range("T3").Value = 0.3172823719910997

# Row 4
$ws.Range("G4").Value = 1.567714666666667
$ws.Range("H4").Value = 4.703144
$ws.Range("I4").Value = 0.4394540882181509
$ws.Range("J4").Value = 0.4394540882181509
$ws.Range("M4").Value = 0.08962966666666666
$ws.Range("N4").Value = 0.268889
$ws.Range("O4").Value = 0.4339761198462219
$ws.Range("P4").Value = 0.4339761198462219
$ws.Range("Q4").Value = 0.1405137430017778
$ws.Range("R4").Value = 1.264623687016
$ws.Range("S4").Value = 0.1907125800554724
$ws.Range("T4").Value = 0.1907125800554724

# Row 5
$ws.Range("G5").Value = 1.567714666666667
$ws.Range("H5").Value = 4.703144
$ws.Range("I5").Value = 0.4394540882181509
$ws.Range("J5").Value = 0.4394540882181509
$ws.Range("M5").Value = 0.1169016666666667
$ws.Range("N5").Value = 0.350705
$ws.Range("O5").Value = 0.5660238801537781
$ws.Range("P5").Value = 0.5660238801537781
$ws.Range("Q5").Value = 0.1832684573911111
$ws.Range("R5").Value = 1.64941611652
$ws.Range("S5").Value = 0.2487415081626785
$ws.Range("T5").Value = 0.2487415081626785
